$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 468.0476
$ws.Range("I33").Value = 459.875
$ws.Range("K33").Value = 459.875
$ws.Range("M33").Value = -230.875
$ws.Range("H100").Value = 1035.4166
$ws.Range("I100").Value = 720
$ws.Range("K100").Value = 720
$ws.Range("M100").Value = -179
$ws.Range("H106").Value = 6859.2085
$ws.Range("I106").Value = 7027
$ws.Range("K106").Value = 7027
$ws.Range("M106").Value = -6396
$ws.Range("H111").Value = 8100
$ws.Range("J111").Value = 7150
$ws.Range("L111").Value = 21450
$ws.Range("N111").Value = -27584
$ws.Range("H138").Value = 440808.03
$ws.Range("I138").Value = 1236.9131
$ws.Range("J138").Value = 573836.1
$ws.Range("K138").Value = 3710.7393
$ws.Range("L138").Value = 1721508.3
$ws.Range("M138").Value = 1429.2607
$ws.Range("N138").Value = -1731788.3
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5722.099
$ws.Range("I32").Value = 4922.383
$ws.Range("J32").Value = 12199.8
$ws.Range("K32").Value = 4922.383
$ws.Range("L32").Value = 12199.8
$ws.Range("M32").Value = -4635.383
$ws.Range("N32").Value = -12773.8
$ws.Range("H45").Value = 1225.5454
$ws.Range("I45").Value = 1148.1
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1148.1
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -771.0999999999999
$ws.Range("N45").Value = -2754
$ws.Range("H132").Value = 2138.6938
$ws.Range("I132").Value = 1789.0333
$ws.Range("K132").Value = 5367.0999
$ws.Range("M132").Value = -2837.0999
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 397.33334
$ws.Range("I22").Value = 193.33333
$ws.Range("J22").Value = 465.33334
$ws.Range("K22").Value = 193.33333
$ws.Range("L22").Value = 465.33334
$ws.Range("M22").Value = -20.33332999999999
$ws.Range("N22").Value = -811.33334
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1643.5778
$ws.Range("J31").Value = 3498.5
$ws.Range("L31").Value = 3498.5
$ws.Range("N31").Value = -4088.5
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H34").Value = 1643.5778
$ws.Range("J34").Value = 3498.5
$ws.Range("L34").Value = 3498.5
$ws.Range("N34").Value = -3902.5
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("K36").Value = 500
$ws.Range("M36").Value = -112
$ws.Range("H40").Value = 500
$ws.Range("I40").Value = 500
$ws.Range("K40").Value = 500
$ws.Range("M40").Value = -340
$ws.Range("H58").Value = 3909.805
$ws.Range("I58").Value = 1160.8462
$ws.Range("J58").Value = 8674.666999999999
$ws.Range("K58").Value = 1160.8462
$ws.Range("L58").Value = 8674.666999999999
$ws.Range("M58").Value = -957.8462
$ws.Range("N58").Value = -9080.666999999999
$ws.Range("H62").Value = 6668858
$ws.Range("I62").Value = 2320.3635
$ws.Range("K62").Value = 2320.3635
$ws.Range("M62").Value = -1696.3635
$ws.Range("H65").Value = 6668858
$ws.Range("I65").Value = 2320.3635
$ws.Range("K65").Value = 11601.8175
$ws.Range("M65").Value = -8481.817499999999
$ws.Range("H107").Value = 1413.7142
$ws.Range("I107").Value = 679.2
$ws.Range("J107").Value = 3250
$ws.Range("K107").Value = 679.2
$ws.Range("L107").Value = 3250
$ws.Range("M107").Value = 1240.8
$ws.Range("N107").Value = -7090
$ws.Range("H132").Value = 2198
$ws.Range("I132").Value = 1853.5555
$ws.Range("J132").Value = 2818
$ws.Range("K132").Value = 5560.666499999999
$ws.Range("L132").Value = 8454
$ws.Range("M132").Value = -3030.666499999999
$ws.Range("N132").Value = -13514
$ws.Range("H136").Value = 3909.805
$ws.Range("I136").Value = 1160.8462
$ws.Range("J136").Value = 8674.666999999999
$ws.Range("K136").Value = 3482.5386
$ws.Range("L136").Value = 26024.001
$ws.Range("M136").Value = -932.5385999999999
$ws.Range("N136").Value = -31124.001
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1316.7894
$ws.Range("J34").Value = 1463.4706
$ws.Range("L34").Value = 4390.4118
$ws.Range("N34").Value = -4558.4118
$ws.Range("H39").Value = 4182.8335
$ws.Range("J39").Value = 4259.4
$ws.Range("L39").Value = 12778.2
$ws.Range("N39").Value = -13366.2
$ws.Range("H55").Value = 2525.4443
$ws.Range("I55").Value = 454
$ws.Range("J55").Value = 2784.375
$ws.Range("K55").Value = 1362
$ws.Range("L55").Value = 8353.125
$ws.Range("M55").Value = -1185
$ws.Range("N55").Value = -8707.125
$ws.Range("H109").Value = 64279.938
$ws.Range("I109").Value = 84423.25
$ws.Range("J109").Value = 3850
$ws.Range("K109").Value = 253269.75
$ws.Range("L109").Value = 11550
$ws.Range("M109").Value = -252229.75
$ws.Range("N109").Value = -13630
$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -460
$ws.Range("N126").ClearContents()
$ws.Range("H129").Value = 26042908
$ws.Range("I129").Value = 83334070
$ws.Range("J129").Value = 6945853.5
$ws.Range("K129").Value = 250002210
$ws.Range("L129").Value = 20837560.5
$ws.Range("M129").Value = -249997210
$ws.Range("N129").Value = -20847560.5
$ws.Range("H131").Value = 14926329
$ws.Range("J131").Value = 997.65
$ws.Range("L131").Value = 2992.95
$ws.Range("N131").Value = -13072.95
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1112.0416
$ws.Range("I113").Value = 1007.3158
$ws.Range("J113").Value = 1510
$ws.Range("K113").Value = 1007.3158
$ws.Range("L113").Value = 1510
$ws.Range("M113").Value = 1162.6842
$ws.Range("N113").Value = -5850
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5198.364
$ws.Range("J46").Value = 6177.778
$ws.Range("L46").Value = 6177.778
$ws.Range("N46").Value = -6553.778
$ws.Range("H61").Value = 1333.3334
$ws.Range("I61").Value = 1333.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1333.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1131.3334
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1333.3334
$ws.Range("I113").Value = 1333.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1333.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 836.6666
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2919.35
$ws.Range("I132").Value = 2820.4546
$ws.Range("K132").Value = 8461.363799999999
$ws.Range("M132").Value = -5931.363799999999
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 370.8
$ws.Range("I107").Value = 312.55554
$ws.Range("K107").Value = 937.66662
$ws.Range("M107").Value = 982.33338
$ws.Range("H113").Value = 315.54544
$ws.Range("I113").Value = 275.83334
$ws.Range("J113").Value = 363.2
$ws.Range("K113").Value = 827.5000200000001
$ws.Range("L113").Value = 1089.6
$ws.Range("M113").Value = 1342.49998
$ws.Range("N113").Value = -5429.6
$ws.Range("H136").Value = 1210.0646
$ws.Range("I136").Value = 1039.92
$ws.Range("K136").Value = 3119.76
$ws.Range("M136").Value = -569.7600000000002
